{"js": "// \"Version 1.\" -> \"Version 2.\"\n//\n// The target revision bumps the version digit and, as a side effect of how\n// Word recorded that edit, leaves \"Version\" split into two runs (\"Versi\" /\n// \"on\") and moves the auto-maintained \"_GoBack\" bookmark so it sits right\n// after the new \"2\" (before the trailing \".\") instead of at the very end of\n// the paragraph. We reproduce both the visible text change and that run /\n// bookmark layout.\n\nconst body = context.document.body;\n\n// Step 1: bump the version number itself, \"1\" -> \"2\".\nconst digit = body.search(\"1\", { matchCase: true, matchWholeWord: false });\ndigit.load(\"text\");\nawait context.sync();\nif (digit.items.length === 0) {\n  throw new Error(\"Could not find the version digit '1' to replace.\");\n}\ndigit.items[0].insertText(\"2\", Word.InsertLocation.replace);\nawait context.sync();\n\n// Step 2: split \"Version\" into \"Versi\" + \"on\" runs. Inserting and\n// immediately removing a bookmark at that boundary forces Word to break the\n// text run there without leaving any stray formatting behind.\nconst prefix = body.search(\"Versi\", { matchCase: true });\nawait context.sync();\nconst splitPoint = prefix.items[0].getRange(\"End\");\nsplitPoint.insertBookmark(\"__tmp_run_split\");\nawait context.sync();\ncontext.document.deleteBookmark(\"__tmp_run_split\");\nawait context.sync();\n\n// Step 3: move the \"_GoBack\" bookmark so it wraps the point right after the\n// newly typed \"2\" (and before the final \".\"), matching where Word leaves it\n// after typing over the old digit.\nconst newDigit = body.search(\"2\", { matchCase: true });\nawait context.sync();\nconst goBackPoint = newDigit.items[0].getRange(\"End\");\ncontext.document.deleteBookmark(\"_GoBack\");\ngoBackPoint.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# \"Version 1.\" -> \"Version 2.\"\n#\n# The target revision bumps the version digit and, as a side effect of how\n# Word recorded that edit, leaves \"Version\" split into two runs (\"Versi\" /\n# \"on\") and moves the auto-maintained \"_GoBack\" bookmark so it sits right\n# after the new \"2\" (before the trailing \".\") instead of at the very end of\n# the paragraph. We reproduce both the visible text change and that run /\n# bookmark layout.\n\n$d = $word.ActiveDocument\n\n# Step 1: bump the version number itself, \"1\" -> \"2\".\n$find = $d.Content\n$find.Find.ClearFormatting()\n$find.Find.Execute(\"1\") | Out-Null\nif (-not $find.Find.Found) {\n    throw \"Could not find the version digit '1' to replace.\"\n}\n$find.Text = \"2\"\n\n# Step 2: split \"Version\" into \"Versi\" + \"on\" runs. Adding and immediately\n# deleting a bookmark at that boundary forces Word to break the text run\n# there without leaving any stray formatting behind.\n$d.Bookmarks.Add(\"__tmp_run_split\", $d.Range(0, 5)) | Out-Null\n$d.Bookmarks(\"__tmp_run_split\").Delete()\n\n# Step 3: move the \"_GoBack\" bookmark so it wraps the point right after the\n# newly typed \"2\" (and before the final \".\"), matching where Word leaves it\n# after typing over the old digit.\n$d.Bookmarks.Add(\"_GoBack\", $d.Range(9, 9)) | Out-Null\n"}
